$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.549.55"
$ws.Range("E2").Value = "  -4.32%  "
$ws.Range("D3").Value = "2.220.34"
$ws.Range("E3").Value = "  -5.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "482.79"
$ws.Range("E5").Value = "  -3.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.07"
$ws.Range("E6").Value = "  -2.47%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -4.89%  "
$ws.Range("D9").Value = "2.230.48"
$ws.Range("E9").Value = "  -5.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0910"
$ws.Range("E10").Value = "  -6.51%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.68"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("D14").Value = "2.614.81"
$ws.Range("E14").Value = "  -5.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.98"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "53.516.88"
$ws.Range("E16").Value = "  -4.37%  "
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").Value = "2.212.35"
$ws.Range("E18").Value = "  -5.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.54"
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.94"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "297.95"
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.09"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.26"
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.996"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.362"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.143"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.95"
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.83"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("E30").Value = "  -3.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "0.0₃0672"
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.71"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.38"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.15"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.829"
$ws.Range("E38").Value = "  +5.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.56"
$ws.Range("E39").Value = "  -4.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.73"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.364"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.36"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.26"
$ws.Range("E43").Value = "  -1.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.33"
$ws.Range("E44").Value = "  -5.09%  "
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.531"
$ws.Range("E47").Value = "  -4.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "229.47"
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0468"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("E51").Value = "  -4.79%  "
